$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: update the "x","y","z" value columns (C:E) with new numbers ---
$ws1.Range("C2").Value = 0.34
$ws1.Range("D2").Value = 8
$ws1.Range("E2").Value = 4

$ws1.Range("C3").Value = 0.316
$ws1.Range("D3").Value = 7.7
$ws1.Range("E3").Value = 5

$ws1.Range("C4").Value = 0.608
$ws1.Range("D4").Value = 7.8
$ws1.Range("E4").Value = 6

$ws1.Range("C5").Value = 0.3
$ws1.Range("D5").Value = 7.9
$ws1.Range("E5").Value = 3

# --- Sheet2: rename class "A2" to "A0" and update the value columns (B:D) ---
$ws2.Range("A3").Value = "A0"

$ws2.Range("B2").Value = 0.3
$ws2.Range("C2").Value = 7.6
$ws2.Range("D2").Value = 2.9

$ws2.Range("B3").Value = 0.7
$ws2.Range("C3").Value = 8
$ws2.Range("D3").Value = 7

# --- View state: make Sheet1 the active/selected tab with B11 selected ---
# and leave Sheet2's selection at A4 (its previous selection is dropped).
$ws2.Activate()
$ws2.Range("A4").Select()

$ws1.Activate()
$ws1.Range("B11").Select()
